$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the year header row (previously 2-digit years, now full 4-digit years).
# Force the cells to remain text (not numeric) by pre-setting a text number format.
$headerRange = $ws.Range("B1:F1")
$headerRange.NumberFormat = "@"

$ws.Range("B1").Value = "2006"
$ws.Range("C1").Value = "2007"
$ws.Range("D1").Value = "2008"
$ws.Range("E1").Value = "2009"
$ws.Range("F1").Value = "2010"
